$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "39.502.17"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.226.00"
$ws.Range("E3").Value = "  -3.96%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'297.69"
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("D6").Value = "'82.61"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("E7").Value = "  -2.47%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.472"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").Value = "'0.0777"
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("D11").Value = "'29.86"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "'46.64"
$ws.Range("E12").Value = "  -11.32%  "
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").Value = "2.570.67"
$ws.Range("E14").Value = "  -3.71%  "
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").Value = "'14.15"
$ws.Range("E16").Value = "  -3.16%  "
$ws.Range("D17").Value = "2.225.53"
$ws.Range("E17").Value = "  -4.03%  "
$ws.Range("D18").Value = "'0.719"
$ws.Range("E18").Value = "  -4.10%  "
$ws.Range("D19").Value = "39.417.22"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").Value = "'5.77"
$ws.Range("E21").Value = "  -4.52%  "
$ws.Range("D22").Value = "'65.09"
$ws.Range("E22").Value = "  -3.34%  "
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").Value = "'229.25"
$ws.Range("E24").Value = "  -2.33%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "  -4.61%  "
$ws.Range("D27").Value = "'1.82"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").Value = "'22.76"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "'32.29"
$ws.Range("E31").Value = "  -6.85%  "
$ws.Range("D32").Value = "'148.78"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Value = "'4.85"
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("E35").Value = "  -2.25%  "
$ws.Range("E36").Value = "  -3.72%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.111"
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "'15.80"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D39").Value = "'0.0970"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("E42").Value = "  -2.38%  "
$ws.Range("D43").Value = "1.925.40"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").Value = "'0.0263"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("E45").Value = "  -9.80%  "
$ws.Range("D46").Value = "'9.22"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "'16.43"
$ws.Range("E47").Value = "  -6.06%  "
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("D49").Value = "2.439.36"
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("D50").Value = "'71.52"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("D51").Value = "'88.45"
$ws.Range("E51").Value = "  -4.05%  "
